$d = $word.ActiveDocument

# --- 1. First bullet: collapse the multi-run "Decrement the lives..." text
#        (incl. spell-check proofErr markers) down to the new single sentence.
$oldText = "Decrement the lives by one every time pac man dies and visually show one less pac-man."
$newText = "Reset level when all pellets are eaten."
$find1 = $d.Content.Find
$find1.ClearFormatting()
$found1 = $find1.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
Write-Output "Replaced first bullet: $found1"

# --- 2. Second bullet: it is currently an empty list paragraph. Add the new
#        sentence, then re-create the _GoBack bookmark right after it (this
#        is where it ends up after the diff, since it's the most recent edit).
$newSentence = "Have screen blink the same way as original before resetting level after winning level."
$p2 = $d.Paragraphs(4)
$r2 = $p2.Range

# Insert the sentence plus a one-character placeholder marker. Placing the
# bookmark exactly at the paragraph-end boundary confuses this host's range
# resolution (it ends up splitting across paragraphs), so we keep a spare
# trailing character in place while we drop the bookmark, then remove it.
$placeholder = "#"
$r2.InsertBefore($newSentence + $placeholder)
Write-Output "Inserted second bullet text"

# The first Find/Replace leaves Word's auto-tracked _GoBack bookmark sitting
# on that edit; remove it so we can place a single _GoBack where it belongs.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
    Write-Output "Removed stray _GoBack"
}

$p2b = $d.Paragraphs(4)
$bmPos = $p2b.Range.Start + $newSentence.Length
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null
Write-Output "Added _GoBack bookmark"

# Remove the placeholder character now that the bookmark is anchored.
$placeholderRange = $d.Range($bmPos, $bmPos + 1)
$placeholderRange.Delete()
Write-Output "Removed placeholder"
